$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The F column holds date serial numbers (formatted mm/dd/yyyy).
# Shift each date in F2:F7 forward by 2 days.
$ws.Range("F2").Value = 44580
$ws.Range("F3").Value = 44579
$ws.Range("F4").Value = 44578
$ws.Range("F5").Value = 44577
$ws.Range("F6").Value = 44576
$ws.Range("F7").Value = 44575
